$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 512.9666999999999
$ws.Range("I80").Value = 265.94446
$ws.Range("K80").Value = 797.83338
$ws.Range("M80").Value = 200.16662
# Row 83
$ws.Range("H83").Value = 512.9666999999999
$ws.Range("I83").Value = 265.94446
$ws.Range("K83").Value = 2393.50014
$ws.Range("M83").Value = 2598.49986
# Row 112
$ws.Range("H112").Value = 1467.303
$ws.Range("J112").Value = 1640.8214
$ws.Range("L112").Value = 4922.4642
$ws.Range("N112").Value = -7138.4642
# Row 132
$ws.Range("H132").Value = 2564.1177
$ws.Range("I132").Value = 2542.2856
$ws.Range("J132").Value = 2666
$ws.Range("K132").Value = 7626.8568
$ws.Range("L132").Value = 7998
$ws.Range("M132").Value = -5096.8568
$ws.Range("N132").Value = -13058
# Row 137
$ws.Range("H137").Value = 2236.8696
$ws.Range("I137").Value = 2465.4
$ws.Range("J137").Value = 1808.375
$ws.Range("K137").Value = 7396.200000000001
$ws.Range("L137").Value = 5425.125
$ws.Range("M137").Value = -4846.200000000001
$ws.Range("N137").Value = -10525.125

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1794.7179
$ws.Range("I2").Value = 1421.9584
$ws.Range("J2").Value = 2391.1333
$ws.Range("K2").Value = 1421.9584
$ws.Range("L2").Value = 2391.1333
$ws.Range("M2").Value = -1308.9584
$ws.Range("N2").Value = -2617.1333
# Row 32
$ws.Range("H32").Value = 2397.9648
$ws.Range("I32").Value = 1096.0613
$ws.Range("J32").Value = 10372.125
$ws.Range("K32").Value = 1096.0613
$ws.Range("L32").Value = 10372.125
$ws.Range("M32").Value = -809.0613000000001
$ws.Range("N32").Value = -10946.125
# Row 45
$ws.Range("H45").Value = 2084.4736
$ws.Range("I45").Value = 1884.0834
$ws.Range("J45").Value = 2428
$ws.Range("K45").Value = 1884.0834
$ws.Range("L45").Value = 2428
$ws.Range("M45").Value = -1507.0834
$ws.Range("N45").Value = -3182
# Row 74
$ws.Range("H74").Value = 3390.4517
$ws.Range("I74").Value = 3200.389
$ws.Range("J74").Value = 3653.6155
$ws.Range("K74").Value = 3200.389
$ws.Range("L74").Value = 3653.6155
$ws.Range("M74").Value = -2326.389
$ws.Range("N74").Value = -5401.6155
# Row 77
$ws.Range("H77").Value = 3390.4517
$ws.Range("I77").Value = 3200.389
$ws.Range("J77").Value = 3653.6155
$ws.Range("K77").Value = 16001.945
$ws.Range("L77").Value = 18268.0775
$ws.Range("M77").Value = -11633.945
$ws.Range("N77").Value = -27004.0775
# Row 88
$ws.Range("H88").Value = 2448.6
$ws.Range("J88").Value = 2532.2856
$ws.Range("L88").Value = 2532.2856
$ws.Range("N88").Value = -3344.2856
# Row 91
$ws.Range("H91").Value = 2448.6
$ws.Range("J91").Value = 2532.2856
$ws.Range("L91").Value = 2532.2856
$ws.Range("N91").Value = -5340.2856
# Row 116
$ws.Range("H116").Value = 1794.7179
$ws.Range("I116").Value = 1421.9584
$ws.Range("J116").Value = 2391.1333
$ws.Range("K116").Value = 1421.9584
$ws.Range("L116").Value = 2391.1333
$ws.Range("M116").Value = 872.0416
$ws.Range("N116").Value = -6979.1333
# Row 122
$ws.Range("H122").Value = 4128.6216
$ws.Range("I122").Value = 4209.8
$ws.Range("J122").Value = 3959.5
$ws.Range("K122").Value = 12629.4
$ws.Range("L122").Value = 11878.5
$ws.Range("M122").Value = -10179.4
$ws.Range("N122").Value = -16778.5
# Row 132
$ws.Range("H132").Value = 13231.739
$ws.Range("I132").Value = 6955.6
$ws.Range("K132").Value = 20866.8
$ws.Range("M132").Value = -18336.8

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1794.7179
$ws.Range("I3").Value = 1421.9584
$ws.Range("J3").Value = 2391.1333
$ws.Range("K3").Value = 1421.9584
$ws.Range("L3").Value = 2391.1333
$ws.Range("M3").Value = -1307.9584
$ws.Range("N3").Value = -2619.1333
# Row 134
$ws.Range("H134").Value = 4006.8262
$ws.Range("I134").Value = 4282.0713
$ws.Range("J134").Value = 3578.6667
$ws.Range("K134").Value = 12846.2139
$ws.Range("L134").Value = 10736.0001
$ws.Range("M134").Value = -10311.2139
$ws.Range("N134").Value = -15806.0001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1085.5883
$ws.Range("I31").Value = 1063.8334
$ws.Range("J31").Value = 1137.8
$ws.Range("K31").Value = 1063.8334
$ws.Range("L31").Value = 1137.8
$ws.Range("M31").Value = -768.8334
$ws.Range("N31").Value = -1727.8
# Row 34
$ws.Range("H34").Value = 1085.5883
$ws.Range("I34").Value = 1063.8334
$ws.Range("J34").Value = 1137.8
$ws.Range("K34").Value = 1063.8334
$ws.Range("L34").Value = 1137.8
$ws.Range("M34").Value = -861.8334
$ws.Range("N34").Value = -1541.8
# Row 99
$ws.Range("H99").Value = 4197.8
$ws.Range("I99").Value = 3928.5173
$ws.Range("J99").Value = 4907.727
$ws.Range("K99").Value = 3928.5173
$ws.Range("L99").Value = 4907.727
$ws.Range("M99").Value = -2430.5173
$ws.Range("N99").Value = -7903.727
# Row 122
$ws.Range("H122").Value = 2341.75
$ws.Range("I122").Value = 2295.5
$ws.Range("J122").Value = 2434.25
$ws.Range("K122").Value = 6886.5
$ws.Range("L122").Value = 7302.75
$ws.Range("M122").Value = -4436.5
$ws.Range("N122").Value = -12202.75
# Row 126
$ws.Range("H126").Value = 4197.8
$ws.Range("I126").Value = 3928.5173
$ws.Range("J126").Value = 4907.727
$ws.Range("K126").Value = 11785.5519
$ws.Range("L126").Value = 14723.181
$ws.Range("M126").Value = -9315.5519
$ws.Range("N126").Value = -19663.181
# Row 132
$ws.Range("H132").Value = 2567
$ws.Range("I132").Value = 2567
$ws.Range("K132").Value = 7701
$ws.Range("M132").Value = -5171
# Row 134
$ws.Range("H134").Value = 2466.8
$ws.Range("I134").Value = 2466.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7400.400000000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4865.400000000001
$ws.Range("N134").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1128.7222
$ws.Range("I131").Value = 682
$ws.Range("J131").Value = 1486.1
$ws.Range("K131").Value = 2046
$ws.Range("L131").Value = 4458.299999999999
$ws.Range("M131").Value = 2994
$ws.Range("N131").Value = -14538.3

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 5380.7
$ws.Range("I126").Value = 7665.6665
$ws.Range("J126").Value = 4401.4287
$ws.Range("K126").Value = 22996.9995
$ws.Range("L126").Value = 13204.2861
$ws.Range("M126").Value = -20526.9995
$ws.Range("N126").Value = -18144.2861
# Row 132
$ws.Range("H132").Value = 3678
$ws.Range("I132").Value = 2856
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 8568
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -6038
$ws.Range("N132").Value = -18560

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1166.8334
$ws.Range("I46").Value = 567
$ws.Range("K46").Value = 567
$ws.Range("M46").Value = -379
# Row 122
$ws.Range("H122").Value = 8843.565000000001
$ws.Range("I122").Value = 7181.125
$ws.Range("J122").Value = 12643.429
$ws.Range("K122").Value = 21543.375
$ws.Range("L122").Value = 37930.287
$ws.Range("M122").Value = -19093.375
$ws.Range("N122").Value = -42830.287

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1082.4375
$ws.Range("I132").Value = 955.3077
$ws.Range("J132").Value = 1633.3334
$ws.Range("K132").Value = 2865.9231
$ws.Range("L132").Value = 4900.0002
$ws.Range("M132").Value = -335.9231
$ws.Range("N132").Value = -9960.0002
